# Auto-generated Excel COM-interop edit script
# Updates market-data derived cells (H..N columns) per scheduled runner refresh
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 3
$ws.Range("H3").Value = 47571.43
$ws.Range("J3").Value = 47571.43
$ws.Range("L3").Value = 47571.43
$ws.Range("N3").Value = -47799.43

# ALC row 40
$ws.Range("H40").Value = 2298.923
$ws.Range("I40").Value = 1648.8
$ws.Range("K40").Value = 1648.8
$ws.Range("M40").Value = -1473.8

# ALC row 70
$ws.Range("H70").Value = 3500
$ws.Range("I70").Value = 3500
$ws.Range("K70").Value = 10500
$ws.Range("M70").Value = -10230

# ALC row 73
$ws.Range("H73").Value = 3500
$ws.Range("I73").Value = 3500
$ws.Range("K73").Value = 10500
$ws.Range("M73").Value = -9564

# ALC row 102
$ws.Range("H102").Value = 47571.43
$ws.Range("J102").Value = 47571.43
$ws.Range("L102").Value = 47571.43
$ws.Range("N102").Value = -54061.43

# ALC row 132
$ws.Range("H132").Value = 1564.3334
$ws.Range("I132").Value = 1346.5
$ws.Range("K132").Value = 4039.5
$ws.Range("M132").Value = -1509.5

# ALC row 137
$ws.Range("H137").Value = 2749.6667
$ws.Range("I137").Value = 3766.3333
$ws.Range("K137").Value = 11298.9999
$ws.Range("M137").Value = -8748.999899999999

# ALC row 138
$ws.Range("H138").Value = 3000
$ws.Range("I138").Value = 2000
$ws.Range("J138").Value = 4000
$ws.Range("K138").Value = 6000
$ws.Range("L138").Value = 12000
$ws.Range("M138").Value = -860
$ws.Range("N138").Value = -22280

$ws = $wb.Worksheets.Item("ARM")
# ARM row 15
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()

# ARM row 45
$ws.Range("H45").Value = 2666.3333
$ws.Range("I45").Value = 2666.3333
$ws.Range("K45").Value = 2666.3333
$ws.Range("M45").Value = -2289.3333

# ARM row 122
$ws.Range("H122").Value = 4212.125
$ws.Range("I122").Value = 4212.125
$ws.Range("K122").Value = 12636.375
$ws.Range("M122").Value = -10186.375

$ws = $wb.Worksheets.Item("BSM")
# BSM row 9
$ws.Range("H9").Value = 74306.336
$ws.Range("J9").Value = 74306.336
$ws.Range("L9").Value = 74306.336
$ws.Range("N9").Value = -74642.336

# BSM row 54
$ws.Range("H54").Value = 30027.666
$ws.Range("I54").Value = 4083
$ws.Range("J54").Value = 43000
$ws.Range("K54").Value = 4083
$ws.Range("L54").Value = 43000
$ws.Range("M54").Value = -3599
$ws.Range("N54").Value = -43968

# BSM row 95
$ws.Range("H95").Value = 59997
$ws.Range("J95").Value = 59997
$ws.Range("L95").Value = 59997
$ws.Range("N95").Value = -65489

# BSM row 105
$ws.Range("H105").Value = 1755.3334
$ws.Range("I105").Value = 1596.7273
$ws.Range("J105").Value = 3500
$ws.Range("K105").Value = 1596.7273
$ws.Range("L105").Value = 3500
$ws.Range("M105").Value = 150.2727
$ws.Range("N105").Value = -6994

# BSM row 107
$ws.Range("H107").Value = 11
$ws.Range("I107").Value = 11
$ws.Range("K107").Value = 11
$ws.Range("M107").Value = 1909

# BSM row 134
$ws.Range("H134").Value = 2285.8
$ws.Range("I134").Value = 2357.25
$ws.Range("J134").Value = 2000
$ws.Range("K134").Value = 7071.75
$ws.Range("L134").Value = 6000
$ws.Range("M134").Value = -4536.75
$ws.Range("N134").Value = -11070

$ws = $wb.Worksheets.Item("CRP")
# CRP row 33
$ws.Range("H33").Value = 25564.285
$ws.Range("I33").Value = 11790.2
$ws.Range("K33").Value = 11790.2
$ws.Range("M33").Value = -11411.2

# CRP row 43
$ws.Range("H43").Value = 32332.334
$ws.Range("J43").Value = 32332.334
$ws.Range("L43").Value = 32332.334
$ws.Range("N43").Value = -32700.334

# CRP row 52
$ws.Range("H52").Value = 75000
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 75000
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 75000
$ws.Range("N52").Value = -75588
$ws.Range("M52").ClearContents()

# CRP row 99
$ws.Range("H99").Value = 2002602.4
$ws.Range("I99").Value = 1253253
$ws.Range("K99").Value = 1253253
$ws.Range("M99").Value = -1251755

# CRP row 101
$ws.Range("H101").Value = 32332.334
$ws.Range("J101").Value = 32332.334
$ws.Range("L101").Value = 32332.334
$ws.Range("N101").Value = -38822.334

# CRP row 126
$ws.Range("H126").Value = 2002602.4
$ws.Range("I126").Value = 1253253
$ws.Range("K126").Value = 3759759
$ws.Range("M126").Value = -3757289

$ws = $wb.Worksheets.Item("CUL")
# CUL row 39
$ws.Range("H39").Value = 16134.5
$ws.Range("J39").Value = 16134.5
$ws.Range("L39").Value = 48403.5
$ws.Range("N39").Value = -48991.5

# CUL row 50
$ws.Range("H50").Value = 1452.25
$ws.Range("I50").Value = 404.5
$ws.Range("K50").Value = 1213.5
$ws.Range("M50").Value = -732.5

# CUL row 53
$ws.Range("H53").Value = 1452.25
$ws.Range("I53").Value = 404.5
$ws.Range("K53").Value = 1213.5
$ws.Range("M53").Value = -732.5

# CUL row 60
$ws.Range("H60").Value = 258.66666
$ws.Range("I60").Value = 138
$ws.Range("K60").Value = 414
$ws.Range("M60").Value = -163

# CUL row 68
$ws.Range("H68").Value = 603.6667
$ws.Range("I68").Value = 595.5
$ws.Range("K68").Value = 1786.5
$ws.Range("M68").Value = -975.5

# CUL row 71
$ws.Range("H71").Value = 603.6667
$ws.Range("I71").Value = 595.5
$ws.Range("K71").Value = 5359.5
$ws.Range("M71").Value = -1303.5

# CUL row 92
$ws.Range("H92").Value = 567.8570999999999
$ws.Range("I92").Value = 596.1667
$ws.Range("J92").Value = 398
$ws.Range("K92").Value = 1788.5001
$ws.Range("L92").Value = 1194
$ws.Range("M92").Value = -540.5001
$ws.Range("N92").Value = -3690

# CUL row 115
$ws.Range("H115").Value = 0
$ws.Range("I115").Value = 0
$ws.Range("K115").Value = 0
$ws.Range("M115").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
# GSM row 122
$ws.Range("H122").Value = 4003.4285
$ws.Range("I122").Value = 2566.5
$ws.Range("K122").Value = 7699.5
$ws.Range("M122").Value = -5249.5

# GSM row 132
$ws.Range("H132").Value = 1097.2
$ws.Range("I132").Value = 1097.2
$ws.Range("K132").Value = 3291.6
$ws.Range("M132").Value = -761.6000000000004

$ws = $wb.Worksheets.Item("LTW")
# LTW row 7
$ws.Range("H7").Value = 1266.3334
$ws.Range("I7").Value = 899.5
$ws.Range("J7").Value = 2000
$ws.Range("K7").Value = 899.5
$ws.Range("L7").Value = 2000
$ws.Range("M7").Value = -787.5
$ws.Range("N7").Value = -2224

# LTW row 16
$ws.Range("H16").Value = 17390.572
$ws.Range("J16").Value = 39966.332
$ws.Range("L16").Value = 39966.332
$ws.Range("N16").Value = -40306.332

# LTW row 68
$ws.Range("H68").Value = 1889.6
$ws.Range("I68").Value = 1889.6
$ws.Range("K68").Value = 1889.6
$ws.Range("M68").Value = -1140.6

# LTW row 71
$ws.Range("H71").Value = 1889.6
$ws.Range("I71").Value = 1889.6
$ws.Range("K71").Value = 9448
$ws.Range("M71").Value = -5704

# LTW row 126
$ws.Range("H126").Value = 1266.3334
$ws.Range("I126").Value = 899.5
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 2698.5
$ws.Range("L126").Value = 6000
$ws.Range("M126").Value = -228.5
$ws.Range("N126").Value = -10940

# LTW row 136
$ws.Range("H136").Value = 2500500
$ws.Range("I136").Value = 2500500
$ws.Range("K136").Value = 7501500
$ws.Range("M136").Value = -7498950

$ws = $wb.Worksheets.Item("WVR")
# WVR row 124
$ws.Range("H124").Value = 49999

Write-Output "edits applied"